# Generate Report for Handoff
# Adds a new row (row 3) to each of the three tables (Overview, zh-cn, de-de)
# describing the hand-off of file
# "2958d5cd-a078-47e6-af8e-3d1068fc7fb7ooo....md"

$wb = $excel.ActiveWorkbook

$baseNameMd   = "2958d5cd-a078-47e6-af8e-3d1068fc7fb7ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$pathAndName  = "e2e\" + $baseNameMd
$ghUrl        = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5912a3559dc7abea230a0d22b3ed398eb199709b/e2e/" + $baseNameMd

$zhXlf = "2958d5cd-a078-47e6-af8e-3d1068fc7fb7oooooooooooooooooooooooooooooooooooooooo.3a14883d6d4bfacaa7447863fa6b42608e8c9271.zh-cn.xlf"
$deXlf = "2958d5cd-a078-47e6-af8e-3d1068fc7fb7oooooooooooooooooooooooooooooooooooooooo.3a14883d6d4bfacaa7447863fa6b42608e8c9271.de-de.xlf"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------
# 1) Overview sheet: add row describing the new file
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rngOverview = $rowOverview.Range

$rngOverview.Cells.Item(1,1).Value = $baseNameMd
$rngOverview.Cells.Item(1,2).Value = $pathAndName
$rngOverview.Cells.Item(1,3).Value = ".md"
$rngOverview.Cells.Item(1,4).Value = ""
$rngOverview.Cells.Item(1,5).Value = "Ready for handoff"
$rngOverview.Cells.Item(1,6).Value = "Ready for handoff"
$rngOverview.Cells.Item(1,7).Value = "2016-08-17 10:26:09"
$rngOverview.Cells.Item(1,7).NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($rngOverview.Cells.Item(1,2), $ghUrl, [Type]::Missing, [Type]::Missing, $pathAndName)

# Widen Priority / zh-cn / de-de columns (autofit-style growth caused by the
# new, longer "Ready for handoff" values)
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

# ---------------------------------------------------------------
# 2) zh-cn sheet: add row describing the new file's zh-cn handoff
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()
$rngZhCn = $rowZhCn.Range

$rngZhCn.Cells.Item(1,1).Value = $baseNameMd
$rngZhCn.Cells.Item(1,2).Value = ".md"
$rngZhCn.Cells.Item(1,3).Value = "Ready for handoff"
$rngZhCn.Cells.Item(1,4).Value = "e2e"
$rngZhCn.Cells.Item(1,5).Value = "ht"
$rngZhCn.Cells.Item(1,6).Value = "'False"
$rngZhCn.Cells.Item(1,7).Value = $zhXlf
$rngZhCn.Cells.Item(1,8).Value = "2016-08-17 10:25:58"
$rngZhCn.Cells.Item(1,8).NumberFormat = $dateFmt
$rngZhCn.Cells.Item(1,11).Value = "0001-01-01 00:00:00"
$rngZhCn.Cells.Item(1,11).NumberFormat = $dateFmt
$rngZhCn.Cells.Item(1,13).Value = "'True"
$rngZhCn.Cells.Item(1,15).Value = "'False"

$wsZhCn.Hyperlinks.Add($rngZhCn.Cells.Item(1,1), $ghUrl, [Type]::Missing, [Type]::Missing, $baseNameMd)

$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332

# ---------------------------------------------------------------
# 3) de-de sheet: add row describing the new file's de-de handoff
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()
$rngDeDe = $rowDeDe.Range

$rngDeDe.Cells.Item(1,1).Value = $baseNameMd
$rngDeDe.Cells.Item(1,2).Value = ".md"
$rngDeDe.Cells.Item(1,3).Value = "Ready for handoff"
$rngDeDe.Cells.Item(1,4).Value = "e2e"
$rngDeDe.Cells.Item(1,5).Value = "ht"
$rngDeDe.Cells.Item(1,6).Value = "'False"
$rngDeDe.Cells.Item(1,7).Value = $deXlf
$rngDeDe.Cells.Item(1,8).Value = "2016-08-17 10:26:09"
$rngDeDe.Cells.Item(1,8).NumberFormat = $dateFmt
$rngDeDe.Cells.Item(1,11).Value = "0001-01-01 00:00:00"
$rngDeDe.Cells.Item(1,11).NumberFormat = $dateFmt
$rngDeDe.Cells.Item(1,13).Value = "'True"
$rngDeDe.Cells.Item(1,15).Value = "'False"

$wsDeDe.Hyperlinks.Add($rngDeDe.Cells.Item(1,1), $ghUrl, [Type]::Missing, [Type]::Missing, $baseNameMd)

$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332

Write-Host "Generate Report for Handoff: done"
